$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 46 (pushes existing rows 46.. down by one, incl. the
# trailing block of hidden rows 280-293 -> 281-294).
$ws.Rows(46).Insert()

# Populate the new row with the "Vehicle Registration State" mapping entry.
$ws.Range("C46").Value = "Vehicle Registration State"
$ws.Range("E46").Value = "MI"
$ws.Range("F46").Value = "wm-req-doc:WarrantModificationRequest/j:ConveyanceRegistration[@structures:id=/wm-req-doc:WarrantModificationRequest/j:ConveyanceRegistrationAssociation/j:ItemRegistration/@structures:ref]/j:JurisdictionNCICLISCode"

# Match the author's row height for the new row.
$ws.Rows(46).RowHeight = 56

# Reflect the view state recorded in the saved workbook: selection sits on
# the newly-entered cell, with the frozen pane scrolled so row 42 is the
# first visible row beneath the header freeze.
$ws.Application.Goto($ws.Range("A42"), $true)
$ws.Range("E46").Select()
